# The deck has two slides (currently at positions 7 and 8) that were saved
# out of order relative to their intended narrative: position 7 shows
# "Current Challenges" (sldId 263) and position 8 shows
# "Roles & Responsibilities" (sldId 262), but the "Roles & Responsibilities"
# slide should come right after the other "Roles & Responsibilities" slide
# (position 6), with "Current Challenges" following it.
#
# Fix the order by moving the slide currently at position 8
# ("Roles & Responsibilities") to position 7, which pushes
# "Current Challenges" down to position 8.

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(8)
$slide.MoveTo(7)
